$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Estado de Cuenta" detail table (rows 16-23) to:
#  - insert a new record for VANESSA ROJAS OLMOS / periodo 1810 at the top (row 16)
#  - shift the existing MARIA IRENE / MARLON VICENTE records down one row
#  - correct previously-zero "Valor Mora" totals and a couple of other amounts
$ws.Range("C16").Value = "45560342"
$ws.Range("D16").Value = "VANESSA ROJAS OLMOS"
$ws.Range("E16").Value = "1810"
$ws.Range("F16").Value = 60000
$ws.Range("G16").Value = 1500000

$ws.Range("C17").Value = "45537049"
$ws.Range("D17").Value = "MARIA IRENE SEGRERA FUENMAYOR"
$ws.Range("E17").Value = "1812"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1000000

$ws.Range("C18").Value = "91529249"
$ws.Range("D18").Value = "MARLON VICENTE BACCA MEDINA"
$ws.Range("E18").Value = "1812"
$ws.Range("F18").Value = 60000
$ws.Range("G18").Value = 1500000

$ws.Range("C19").Value = "45560342"
$ws.Range("D19").Value = "VANESSA ROJAS OLMOS"
$ws.Range("E19").Value = "1812"
$ws.Range("F19").Value = 60000
$ws.Range("G19").Value = 1500000

$ws.Range("E20").Value = "1812"
$ws.Range("G20").Value = 781242

$ws.Range("E21").Value = "1812"

$ws.Range("E22").Value = "1812"
$ws.Range("G22").Value = 781242

$ws.Range("E23").Value = "1812"
